$wb = $excel.ActiveWorkbook

# Update sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 329
$ws1.Range("F4").Value = 1336

# Update sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 329
$ws4.Range("F4").Value = 1336
